$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '44.842.91'
    'E2' = '  +4.38%  '
    'D3' = '2.424.53'
    'E3' = '  +2.75%  '
    'E4' = '  -0.11%  '
    'D5' = '316.20'
    'E5' = '  +4.38%  '
    'D6' = '101.88'
    'E7' = '  +2.80%  '
    'E8' = '  -0.12%  '
    'D9' = '0.533'
    'E9' = '  +12.08%  '
    'D10' = '35.45'
    'E10' = '  +3.38%  '
    'E11' = '  +1.93%  '
    'D12' = '18.76'
    'E12' = '  +1.48%  '
    'E13' = '  -1.67%  '
    'D14' = '6.93'
    'E14' = '  +3.61%  '
    'D15' = '2.799.84'
    'E15' = '  +2.67%  '
    'D16' = '2.443.90'
    'E16' = '  +2.67%  '
    'E17' = '  +4.73%  '
    'D18' = '44.642.71'
    'E18' = '  +3.97%  '
    'D19' = '12.34'
    'E19' = '  +3.43%  '
    'E20' = '  +2.08%  '
    'D21' = '0.0₃0918'
    'E21' = '  +3.82%  '
    'D22' = '68.81'
    'E22' = '  +1.34%  '
    'D23' = '242.43'
    'E23' = '  +3.14%  '
    'E24' = '  +4.43%  '
    'D25' = '2.50'
    'E25' = '  +2.37%  '
    'E26' = '  -0.02%  '
    'D27' = '25.22'
    'E27' = '  +3.56%  '
    'D28' = '2.19'
    'E28' = '  -7.45%  '
    'D29' = '9.49'
    'E29' = '  +1.73%  '
    'D30' = '33.75'
    'E30' = '  +4.60%  '
    'E31' = '  +1.64%  '
    'E32' = '  +18.58%  '
    'D33' = '19.58'
    'E33' = '  +11.91%  '
    'E34' = '  +3.32%  '
    'D35' = '0.0770'
    'E35' = '  +6.25%  '
    'E36' = '  +0.22%  '
    'E37' = '  +3.91%  '
    'E38' = '  +3.79%  '
    'D39' = '2.86'
    'E39' = '  +0.88%  '
    'D40' = '122.51'
    'E40' = '  -4.82%  '
    'D41' = '0.110'
    'E41' = '  +1.96%  '
    'E42' = '  -3.01%  '
    'E43' = '  +0.35%  '
    'E44' = '  +4.31%  '
    'D45' = '1.943.82'
    'E45' = '  +0.90%  '
    'B46' = 'NEARProtocol'
    'C46' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D46' = '2.93'
    'E46' = '  +8.22%  '
    'B47' = 'ApeXProtocol'
    'C47' = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    'D47' = '2.09'
    'E47' = '  -2.27%  '
    'D48' = '9.32'
    'E48' = '  +0.76%  '
    'D49' = '1.74'
    'E49' = '  +15.50%  '
    'D50' = '75.81'
    'E50' = '  +6.17%  '
    'D51' = '53.87'
    'E51' = '  +5.17%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
